$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, reusing the same formatting as the
# other header cells (e.g. G1) by copying the format over.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the Save values for rows 2-8
$values = @(1, 0, 0, 0, 0, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
